# The commit swaps the deck's colour theme: the slide master (theme1.xml)
# moves from the "Integral" / "Red Violet" palette to the stock "Office
# Theme" / "Office" palette (and, in the canonical package, the old
# "Integral" colours end up preserved in theme2.xml, which backs the
# Notes Master and isn't reachable from slide-level content).
#
# Font scheme and format scheme (fills/lines/effects) are already
# byte-identical between the two themes, so only the 12 theme colours
# need to change. Drive that through the slide's ThemeColorScheme,
# which maps 1:1 onto <a:clrScheme> (dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink) for the presentation's single slide master.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# New "Office Theme" colours, in clrScheme order:
# dk1, lt1, dk2, lt2, accent1, accent2, accent3, accent4, accent5, accent6, hlink, folHlink
$officeThemeColors = @(
    0x000000,  # dk1      000000
    0xFFFFFF,  # lt1      FFFFFF
    0x44546A,  # dk2      44546A
    0xE7E6E6,  # lt2      E7E6E6
    0x5B9BD5,  # accent1  5B9BD5
    0xED7D31,  # accent2  ED7D31
    0xA5A5A5,  # accent3  A5A5A5
    0xFFC000,  # accent4  FFC000
    0x4472C4,  # accent5  4472C4
    0x70AD47,  # accent6  70AD47
    0x0563C1,  # hlink    0563C1
    0x954F72   # folHlink 954F72
)

for ($i = 1; $i -le $tcs.Count; $i++) {
    # RGBColor.RGB takes 0xBBGGRR-ordered (Windows COLORREF) integers, so
    # convert each 0xRRGGBB constant above before assigning.
    $rgbHex = $officeThemeColors[$i - 1]
    $r = ($rgbHex -shr 16) -band 0xFF
    $g = ($rgbHex -shr 8) -band 0xFF
    $b = $rgbHex -band 0xFF
    $colorRef = $r + ($g * 256) + ($b * 65536)

    $tcs.Colors($i).RGB = $colorRef
}
